$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1761006289308176
$ws.Range("C2").Value = 0.5911949685534591
$ws.Range("J2").Value = 0.01886792452830189
$ws.Range("P2").Value = 0.1415094339622641
$ws.Range("S2").Value = 0.07232704402515723

# Row 3
$ws.Range("B3").Value = 0.01025641025641026
$ws.Range("C3").Value = 0.01538461538461539
$ws.Range("J3").Value = 0.03589743589743589
$ws.Range("P3").Value = 0.7948717948717948
$ws.Range("S3").Value = 0.1435897435897436

# Row 4
$ws.Range("J4").Value = 0.08928571428571429
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.1607142857142857

# Row 6
$ws.Range("B6").Value = 0.07531380753138076
$ws.Range("D6").Value = 0.02092050209205021
$ws.Range("F6").Value = 0.07112970711297072
$ws.Range("J6").Value = 0.2426778242677824
$ws.Range("O6").Value = 0.02928870292887029
$ws.Range("Q6").Value = 0.2217573221757322
$ws.Range("R6").Value = 0.03765690376569038
$ws.Range("S6").Value = 0.301255230125523

# Row 7
$ws.Range("B7").Value = 0.08866995073891626
$ws.Range("D7").Value = 0.01970443349753695
$ws.Range("F7").Value = 0.0541871921182266
$ws.Range("J7").Value = 0.187192118226601
$ws.Range("O7").Value = 0.01477832512315271
$ws.Range("Q7").Value = 0.2068965517241379
$ws.Range("R7").Value = 0.06896551724137931
$ws.Range("S7").Value = 0.3596059113300493

# Row 8
$ws.Range("B8").Value = 0.09381663113006397
$ws.Range("D8").Value = 0.02771855010660981
$ws.Range("E8").Value = 0.002132196162046908
$ws.Range("F8").Value = 0.05970149253731343
$ws.Range("J8").Value = 0.09594882729211088
$ws.Range("O8").Value = 0.03198294243070363
$ws.Range("Q8").Value = 0.2281449893390192
$ws.Range("R8").Value = 0.07889125799573561
$ws.Range("S8").Value = 0.3816631130063966

# Row 9
$ws.Range("B9").Value = 0.10625
$ws.Range("D9").Value = 0.025
$ws.Range("F9").Value = 0.05
$ws.Range("O9").Value = 0.03125
$ws.Range("Q9").Value = 0.18125
$ws.Range("R9").Value = 0.05
$ws.Range("S9").Value = 0.43125

# Row 10
$ws.Range("B10").Value = 0.1144662921348315
$ws.Range("D10").Value = 0.02176966292134832
$ws.Range("E10").Value = 0.002106741573033708
$ws.Range("F10").Value = 0.07654494382022473
$ws.Range("J10").Value = 0.113061797752809
$ws.Range("O10").Value = 0.01825842696629213
$ws.Range("Q10").Value = 0.2443820224719101
$ws.Range("R10").Value = 0.06460674157303371
$ws.Range("S10").Value = 0.3448033707865168

# Row 11
$ws.Range("G11").Value = 0.1073619631901841
$ws.Range("J11").Value = 0.09202453987730061
$ws.Range("K11").Value = 0.1779141104294479
$ws.Range("L11").Value = 0.6012269938650306
$ws.Range("S11").Value = 0.02147239263803681

# Row 12
$ws.Range("G12").Value = 0.7277227722772277
$ws.Range("J12").Value = 0.2326732673267327
$ws.Range("K12").Value = 0.004950495049504951
$ws.Range("L12").Value = 0.0198019801980198
$ws.Range("S12").Value = 0.01485148514851485

# Row 13
$ws.Range("F13").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.6046511627906976
$ws.Range("J13").Value = 0.3255813953488372
$ws.Range("S13").Value = 0.04651162790697674

# Row 15
$ws.Range("F15").Value = 0.01879699248120301
$ws.Range("H15").Value = 0.1729323308270677
$ws.Range("I15").Value = 0.06015037593984962
$ws.Range("J15").Value = 0.3120300751879699
$ws.Range("K15").Value = 0.07894736842105263
$ws.Range("M15").Value = 0.01879699248120301
$ws.Range("O15").Value = 0.04135338345864661
$ws.Range("S15").Value = 0.2969924812030075

# Row 16
$ws.Range("F16").Value = 0.02521008403361345
$ws.Range("H16").Value = 0.1512605042016807
$ws.Range("I16").Value = 0.05042016806722689
$ws.Range("J16").Value = 0.4579831932773109
$ws.Range("K16").Value = 0.1050420168067227
$ws.Range("M16").Value = 0.01680672268907563
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.134453781512605

# Row 17
$ws.Range("F17").Value = 0.02413793103448276
$ws.Range("H17").Value = 0.1844827586206897
$ws.Range("I17").Value = 0.0603448275862069
$ws.Range("J17").Value = 0.4448275862068966
$ws.Range("K17").Value = 0.09655172413793103
$ws.Range("M17").Value = 0.02413793103448276
$ws.Range("O17").Value = 0.06379310344827586
$ws.Range("S17").Value = 0.1017241379310345

# Row 18
$ws.Range("F18").Value = 0.02484472049689441
$ws.Range("H18").Value = 0.1801242236024845
$ws.Range("I18").Value = 0.09937888198757763
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.09316770186335403
$ws.Range("M18").Value = 0.0124223602484472
$ws.Range("O18").Value = 0.07453416149068323
$ws.Range("S18").Value = 0.08695652173913043

# Row 19
$ws.Range("F19").Value = 0.01348136399682792
$ws.Range("H19").Value = 0.1982553528945281
$ws.Range("I19").Value = 0.06582077716098335
$ws.Range("J19").Value = 0.3965107057890563
$ws.Range("K19").Value = 0.1165741475019826
$ws.Range("M19").Value = 0.01348136399682792
$ws.Range("N19").Value = 0.0007930214115781126
$ws.Range("O19").Value = 0.08247422680412371
$ws.Range("S19").Value = 0.112609040444092
